$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44187
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/caja 15 kilos granel"
$ws.Range("S2").Value = 1000

# Row 3
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 95
$ws.Range("N3").Value = 13500
$ws.Range("O3").Value = 13500
$ws.Range("P3").Value = 13500
$ws.Range("S3").Value = 900

# Row 4
$ws.Range("L4").Value = "Tercera"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("S4").Value = 800

# Row 5
$ws.Range("D5").Value = 44160
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 44160
$ws.Range("K6").Value = "Castle Brite"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 7000
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S6").Value = 700

# Row 7
$ws.Range("D7").Value = 44159
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S7").Value = 800
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44159
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S8").Value = 700
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44175
$ws.Range("K9").Value = "Modesto"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 140
$ws.Range("N9").Value = 11000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 11571
$ws.Range("Q9").Value = "$/caja 12 kilos"
$ws.Range("S9").Value = 964
$ws.Range("T9").Value = 12

# Row 10
$ws.Range("D10").Value = 44162
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 8500
$ws.Range("O10").Value = 8500
$ws.Range("P10").Value = 8500
$ws.Range("S10").Value = 850

# Row 11
$ws.Range("D11").Value = 44162
$ws.Range("M11").Value = 75
$ws.Range("N11").Value = 14000
$ws.Range("P11").Value = 14400
$ws.Range("S11").Value = 800

# Row 12
$ws.Range("D12").Value = 44167
$ws.Range("K12").Value = "Castle Brite"
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 85
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 10

# Row 13
$ws.Range("D13").Value = 44167
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 9500
$ws.Range("O13").Value = 9500
$ws.Range("P13").Value = 9500
$ws.Range("S13").Value = 950

# Row 14
$ws.Range("D14").Value = 44167
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 833

# Row 15
$ws.Range("D15").Value = 44174
$ws.Range("K15").Value = "Modesto"
$ws.Range("N15").Value = 8500
$ws.Range("O15").Value = 8500
$ws.Range("P15").Value = 8500
$ws.Range("Q15").Value = "$/bandeja 10 kilos"
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 850
$ws.Range("T15").Value = 10

# Row 16
$ws.Range("D16").Value = 44174
$ws.Range("K16").Value = "Modesto"
$ws.Range("M16").Value = 180
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("Q16").Value = "$/caja 18 kilos"
$ws.Range("R16").Value = "Región Metropolitana"
$ws.Range("S16").Value = 833
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44174
$ws.Range("K17").Value = "Modesto"
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("S17").Value = 667
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44189
$ws.Range("K18").Value = "Patterson"
$ws.Range("M18").Value = 130
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 12000
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S18").Value = 667
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44169
$ws.Range("K19").Value = "Dina"
$ws.Range("M19").Value = 80
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 10000
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1000

# Row 20
$ws.Range("D20").Value = 44195
$ws.Range("K20").Value = "Patterson"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 124
$ws.Range("N20").Value = 13000
$ws.Range("O20").Value = 13000
$ws.Range("P20").Value = 13000
$ws.Range("Q20").Value = "$/caja 15 kilos"
$ws.Range("S20").Value = 867
$ws.Range("T20").Value = 15

# Row 21
$ws.Range("D21").Value = 44176
$ws.Range("K21").Value = "Modesto"
$ws.Range("M21").Value = 115
$ws.Range("N21").Value = 11000
$ws.Range("P21").Value = 11609
$ws.Range("Q21").Value = "$/caja 12 kilos"
$ws.Range("S21").Value = 967
$ws.Range("T21").Value = 12

# Row 24
$ws.Range("D24").Value = 44194
$ws.Range("K24").Value = "Patterson"
$ws.Range("M24").Value = 120
$ws.Range("N24").Value = 13000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 13000
$ws.Range("Q24").Value = "$/caja 15 kilos"
$ws.Range("S24").Value = 867
$ws.Range("T24").Value = 15
